# Reorders the detail rows (16-19) of the account-statement table.
# The block of 4 worker/period rows gets reversed:
#   row16 <-> row19, row17 <-> row18
# (columns C:G - "N Doc Trabajador", "Nombre Trabajador", "Periodo Mora",
#  "Valor Mora", "Salario Basico"). Column B ("Tipo Doc Trabajador" = CC)
# and row 20 are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 19
$firstCol = 3   # C
$lastCol  = 7   # G

# Snapshot current values for the block before overwriting anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the snapshot back in reverse row order.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $lastRow - ($r - $firstRow)
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$c]
    }
}
